# Append new daily data rows (2024-09-30 through 2025-01-14) to the
# "nifty smallcap 250" worksheet, extending the table from row 3408 to row 3481.
# Columns: A=Date (serial), B=Total Returns Index, C=P/E, D=P/B, E=Div Yield

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 3409
$nrows = 73
$endRow = $startRow + $nrows - 1

# Build a 2-D array with all the new row values (date serial + 4 numeric columns)
$arr = New-Object 'object[,]' $nrows,5

$arr[0,0] = 45565; $arr[0,1] = 23256.45; $arr[0,2] = 33.49; $arr[0,3] = 4.32; $arr[0,4] = 0.83
$arr[1,0] = 45566; $arr[1,1] = 23409.41; $arr[1,2] = 33.71; $arr[1,3] = 4.35; $arr[1,4] = 0.83
$arr[2,0] = 45568; $arr[2,1] = 22964.06; $arr[2,2] = 33.07; $arr[2,3] = 4.27; $arr[2,4] = 0.84
$arr[3,0] = 45569; $arr[3,1] = 22765.03; $arr[3,2] = 32.78; $arr[3,3] = 4.23; $arr[3,4] = 0.85
$arr[4,0] = 45572; $arr[4,1] = 22124.55; $arr[4,2] = 31.86; $arr[4,3] = 4.11; $arr[4,4] = 0.88
$arr[5,0] = 45573; $arr[5,1] = 22589.83; $arr[5,2] = 32.53; $arr[5,3] = 4.2; $arr[5,4] = 0.86
$arr[6,0] = 45574; $arr[6,1] = 22869.11; $arr[6,2] = 32.93; $arr[6,3] = 4.25; $arr[6,4] = 0.85
$arr[7,0] = 45575; $arr[7,1] = 22898.67; $arr[7,2] = 33.15; $arr[7,3] = 4.32; $arr[7,4] = 0.85
$arr[8,0] = 45576; $arr[8,1] = 23006.92; $arr[8,2] = 33.3; $arr[8,3] = 4.34; $arr[8,4] = 0.84
$arr[9,0] = 45579; $arr[9,1] = 23072.73; $arr[9,2] = 33.4; $arr[9,3] = 4.35; $arr[9,4] = 0.84
$arr[10,0] = 45580; $arr[10,1] = 23328.56; $arr[10,2] = 33.77; $arr[10,3] = 4.4; $arr[10,4] = 0.83
$arr[11,0] = 45581; $arr[11,1] = 23376.96; $arr[11,2] = 33.85; $arr[11,3] = 4.43; $arr[11,4] = 0.83
$arr[12,0] = 45582; $arr[12,1] = 23046.43; $arr[12,2] = 33.39; $arr[12,3] = 4.36; $arr[12,4] = 0.84
$arr[13,0] = 45583; $arr[13,1] = 23037.1; $arr[13,2] = 33.37; $arr[13,3] = 4.36; $arr[13,4] = 0.84
$arr[14,0] = 45586; $arr[14,1] = 22657.48; $arr[14,2] = 32.79; $arr[14,3] = 4.29; $arr[14,4] = 0.85
$arr[15,0] = 45587; $arr[15,1] = 21823.05; $arr[15,2] = 31.37; $arr[15,3] = 4.13; $arr[15,4] = 0.88
$arr[16,0] = 45588; $arr[16,1] = 22051.01; $arr[16,2] = 31.65; $arr[16,3] = 4.18; $arr[16,4] = 0.87
$arr[17,0] = 45589; $arr[17,1] = 21919.67; $arr[17,2] = 31.72; $arr[17,3] = 4.15; $arr[17,4] = 0.88
$arr[18,0] = 45590; $arr[18,1] = 21419.26; $arr[18,2] = 31.01; $arr[18,3] = 4.06; $arr[18,4] = 0.9
$arr[19,0] = 45593; $arr[19,1] = 21666.28; $arr[19,2] = 31.58; $arr[19,3] = 4.1; $arr[19,4] = 0.89
$arr[20,0] = 45594; $arr[20,1] = 21826.31; $arr[20,2] = 31.75; $arr[20,3] = 4.13; $arr[20,4] = 0.88
$arr[21,0] = 45595; $arr[21,1] = 22085.49; $arr[21,2] = 32.13; $arr[21,3] = 4.18; $arr[21,4] = 0.87
$arr[22,0] = 45596; $arr[22,1] = 22424.5; $arr[22,2] = 32.37; $arr[22,3] = 4.25; $arr[22,4] = 0.86
$arr[23,0] = 45597; $arr[23,1] = 22661.58; $arr[23,2] = 32.71; $arr[23,3] = 4.29; $arr[23,4] = 0.84
$arr[24,0] = 45600; $arr[24,1] = 22245.98; $arr[24,2] = 32.03; $arr[24,3] = 4.21; $arr[24,4] = 0.86
$arr[25,0] = 45601; $arr[25,1] = 22353.8; $arr[25,2] = 32.19; $arr[25,3] = 4.23; $arr[25,4] = 0.85
$arr[26,0] = 45602; $arr[26,1] = 22775.3; $arr[26,2] = 32.81; $arr[26,3] = 4.31; $arr[26,4] = 0.84
$arr[27,0] = 45603; $arr[27,1] = 22634.01; $arr[27,2] = 32.64; $arr[27,3] = 4.29; $arr[27,4] = 0.85
$arr[28,0] = 45604; $arr[28,1] = 22282.26; $arr[28,2] = 32.14; $arr[28,3] = 4.22; $arr[28,4] = 0.87
$arr[29,0] = 45607; $arr[29,1] = 22016.49; $arr[29,2] = 31.82; $arr[29,3] = 4.17; $arr[29,4] = 0.87
$arr[30,0] = 45608; $arr[30,1] = 21729.02; $arr[30,2] = 31.47; $arr[30,3] = 4.11; $arr[30,4] = 0.88
$arr[31,0] = 45609; $arr[31,1] = 21064.28; $arr[31,2] = 30.62; $arr[31,3] = 3.99; $arr[31,4] = 0.91
$arr[32,0] = 45610; $arr[32,1] = 21231.61; $arr[32,2] = 30.85; $arr[32,3] = 4.02; $arr[32,4] = 0.91
$arr[33,0] = 45614; $arr[33,1] = 21095.49; $arr[33,2] = 30.54; $arr[33,3] = 3.99; $arr[33,4] = 0.91
$arr[34,0] = 45615; $arr[34,1] = 21289.19; $arr[34,2] = 30.49; $arr[34,3] = 4.03; $arr[34,4] = 0.92
$arr[35,0] = 45617; $arr[35,1] = 21166.06; $arr[35,2] = 32.15; $arr[35,3] = 4.01; $arr[35,4] = 0.91
$arr[36,0] = 45618; $arr[36,1] = 21373.93; $arr[36,2] = 32.52; $arr[36,3] = 4.05; $arr[36,4] = 0.9
$arr[37,0] = 45621; $arr[37,1] = 21825.83; $arr[37,2] = 33.2; $arr[37,3] = 4.13; $arr[37,4] = 0.88
$arr[38,0] = 45622; $arr[38,1] = 21960.17; $arr[38,2] = 33.41; $arr[38,3] = 4.16; $arr[38,4] = 0.87
$arr[39,0] = 45623; $arr[39,1] = 22184.43; $arr[39,2] = 33.75; $arr[39,3] = 4.15; $arr[39,4] = 0.87
$arr[40,0] = 45624; $arr[40,1] = 22216.29; $arr[40,2] = 33.8; $arr[40,3] = 4.09; $arr[40,4] = 0.87
$arr[41,0] = 45625; $arr[41,1] = 22396.2; $arr[41,2] = 34.11; $arr[41,3] = 4.13; $arr[41,4] = 0.88
$arr[42,0] = 45628; $arr[42,1] = 22552.51; $arr[42,2] = 34.34; $arr[42,3] = 4.16; $arr[42,4] = 0.87
$arr[43,0] = 45629; $arr[43,1] = 22764.13; $arr[43,2] = 34.67; $arr[43,3] = 4.2; $arr[43,4] = 0.86
$arr[44,0] = 45630; $arr[44,1] = 22936.66; $arr[44,2] = 34.93; $arr[44,3] = 4.23; $arr[44,4] = 0.86
$arr[45,0] = 45631; $arr[45,1] = 23046.47; $arr[45,2] = 35.1; $arr[45,3] = 4.25; $arr[45,4] = 0.86
$arr[46,0] = 45632; $arr[46,1] = 23171.08; $arr[46,2] = 35.29; $arr[46,3] = 4.27; $arr[46,4] = 0.85
$arr[47,0] = 45635; $arr[47,1] = 23256.39; $arr[47,2] = 35.42; $arr[47,3] = 4.29; $arr[47,4] = 0.85
$arr[48,0] = 45636; $arr[48,1] = 23324.82; $arr[48,2] = 35.52; $arr[48,3] = 4.3; $arr[48,4] = 0.85
$arr[49,0] = 45637; $arr[49,1] = 23390.65; $arr[49,2] = 35.62; $arr[49,3] = 4.31; $arr[49,4] = 0.84
$arr[50,0] = 45638; $arr[50,1] = 23172.77; $arr[50,2] = 35.29; $arr[50,3] = 4.27; $arr[50,4] = 0.85
$arr[51,0] = 45639; $arr[51,1] = 23088.61; $arr[51,2] = 35.16; $arr[51,3] = 4.26; $arr[51,4] = 0.85
$arr[52,0] = 45642; $arr[52,1] = 23190.08; $arr[52,2] = 35.31; $arr[52,3] = 4.28; $arr[52,4] = 0.85
$arr[53,0] = 45643; $arr[53,1] = 23061.18; $arr[53,2] = 35.12; $arr[53,3] = 4.25; $arr[53,4] = 0.85
$arr[54,0] = 45644; $arr[54,1] = 22905.13; $arr[54,2] = 34.88; $arr[54,3] = 4.22; $arr[54,4] = 0.86
$arr[55,0] = 45645; $arr[55,1] = 22827.49; $arr[55,2] = 34.76; $arr[55,3] = 4.21; $arr[55,4] = 0.86
$arr[56,0] = 45646; $arr[56,1] = 22367.9; $arr[56,2] = 34.06; $arr[56,3] = 4.13; $arr[56,4] = 0.88
$arr[57,0] = 45649; $arr[57,1] = 22272.43; $arr[57,2] = 33.92; $arr[57,3] = 4.11; $arr[57,4] = 0.88
$arr[58,0] = 45650; $arr[58,1] = 22341.14; $arr[58,2] = 34.02; $arr[58,3] = 4.03; $arr[58,4] = 0.88
$arr[59,0] = 45652; $arr[59,1] = 22334.78; $arr[59,2] = 34.01; $arr[59,3] = 4.03; $arr[59,4] = 0.88
$arr[60,0] = 45653; $arr[60,1] = 22395.19; $arr[60,2] = 34.1; $arr[60,3] = 4.04; $arr[60,4] = 0.88
$arr[61,0] = 45656; $arr[61,1] = 22297.75; $arr[61,2] = 33.96; $arr[61,3] = 4.03; $arr[61,4] = 0.88
$arr[62,0] = 45657; $arr[62,1] = 22442.03; $arr[62,2] = 34.39; $arr[62,3] = 4.1; $arr[62,4] = 0.87
$arr[63,0] = 45658; $arr[63,1] = 22639.89; $arr[63,2] = 34.7; $arr[63,3] = 4.14; $arr[63,4] = 0.86
$arr[64,0] = 45659; $arr[64,1] = 22765.23; $arr[64,2] = 34.89; $arr[64,3] = 4.16; $arr[64,4] = 0.85
$arr[65,0] = 45660; $arr[65,1] = 22721.16; $arr[65,2] = 34.82; $arr[65,3] = 4.15; $arr[65,4] = 0.85
$arr[66,0] = 45663; $arr[66,1] = 22067.12; $arr[66,2] = 33.82; $arr[66,3] = 4.03; $arr[66,4] = 0.88
$arr[67,0] = 45664; $arr[67,1] = 22379.17; $arr[67,2] = 34.3; $arr[67,3] = 4.09; $arr[67,4] = 0.87
$arr[68,0] = 45665; $arr[68,1] = 22069.9; $arr[68,2] = 33.82; $arr[68,3] = 4.03; $arr[68,4] = 0.88
$arr[69,0] = 45666; $arr[69,1] = 21812.55; $arr[69,2] = 33.43; $arr[69,3] = 3.98; $arr[69,4] = 0.89
$arr[70,0] = 45667; $arr[70,1] = 21280.42; $arr[70,2] = 32.61; $arr[70,3] = 3.89; $arr[70,4] = 0.91
$arr[71,0] = 45670; $arr[71,1] = 20429.01; $arr[71,2] = 31.31; $arr[71,3] = 3.73; $arr[71,4] = 0.95
$arr[72,0] = 45671; $arr[72,1] = 20761.31; $arr[72,2] = 31.84; $arr[72,3] = 3.79; $arr[72,4] = 0.9399999999999999

# Write the whole block in one shot
$rng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 5))
$rng.Value2 = $arr

# Apply the same date number format used by the existing column-A date cells
# (this reuses the existing style index, matching the source workbook's s="5")
$dateColRng = $ws.Range($ws.Cells.Item($startRow, 1), $ws.Cells.Item($endRow, 1))
$dateColRng.NumberFormat = $ws.Cells.Item($startRow - 1, 1).NumberFormat

Write-Host "New UsedRange:" $ws.UsedRange.Address()
